# Append seven new talks (rows 6-12) to the talks worksheet and update
# the active selection, mirroring the "Updated talks and talk map" commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: LIGO Detector Characterization Update ---------------------
$ws.Range("A6").Value = "LIGO Detector Characterization Update"
$ws.Range("B6").Value = "Plenary talk"
$ws.Range("C6").Value = "LVC"
$ws.Range("D6").Value = "LIGO-Virgo Collaboration Meeting"
$ws.Range("E6").Value = "2016-03-17"
$ws.Range("F6").Value = "Pasadena, CA"

# --- Row 7: Burst Data Quality (LVC) -----------------------------------
$ws.Range("A7").Value = "Burst Data Quality "
$ws.Range("C7").Value = "LVC"
$ws.Range("D7").Value = "LIGO-Virgo Collaboration Meeting"
$ws.Range("E7").Value = "2015-09-03"
$ws.Range("F7").Value = "Budapest, Hungary"

# --- Row 8: Intro to Advanced LIGO Burst DQ @ Albert Einstein Institute
$ws.Range("A8").Value = "Introduction to Advanced LIGO Burst Data Quality"
$ws.Range("C8").Value = "seminar"
$ws.Range("D8").Value = "Albert Einstein Institute"
$ws.Range("E8").Value = "2015-08-18"
$ws.Range("F8").Value = "Hannover, Germany"

# --- Row 9: Intro to Advanced LIGO Burst DQ @ Cardiff University ------
$ws.Range("A9").Value = "Introduction to Advanced LIGO Burst Data Quality"
$ws.Range("C9").Value = "seminar"
$ws.Range("D9").Value = "Cardiff University"
$ws.Range("E9").Value = "2015-08-25"
$ws.Range("F9").Value = "Cardiff, Wales"

# --- Row 10: Characterization of the Instrumental Background... -------
$ws.Range("A10").Value = "Characterization of the Instrumental Background of Advanced LIGO's Gravitational Wave Burst Search"
$ws.Range("C10").Value = "GCG"
$ws.Range("D10").Value = "Gulf Coast Gravity Meeting"
$ws.Range("E10").Value = "2015-02-15"
$ws.Range("F10").Value = "Gainesville, FL"

# --- Row 11: Burst Data Quality in Livingston Full Interferometer -----
$ws.Range("A11").Value = "Burst Data Quality in Livingston Full Interferometer"
$ws.Range("C11").Value = "LVC"
$ws.Range("D11").Value = "LIGO-Virgo Collaboration Meeting"
$ws.Range("E11").Value = "2014-08-15"
$ws.Range("F11").Value = "Stanford, CA"

# --- Row 12: Update on Ongoing Projects in Advanced LIGO Suspensions ---
$ws.Range("A12").Value = "Update on Ongoing Projects in Advanced LIGO Suspensions Detector Characterization"
$ws.Range("C12").Value = "LVC"
$ws.Range("D12").Value = "LIGO-Virgo Collaboration Meeting"
$ws.Range("E12").Value = "2013-09-17"
$ws.Range("F12").Value = "Hannover, Germany"

# --- Formatting: mirror the existing look (column A / E / some D cells) -
# Column A cells use the same "s=2" (Arial Unicode MS 10pt) style as the
# rest of the title column.
$ws.Range("A2").Copy()
$ws.Range("A6:A12").PasteSpecial(-4122)

# Column E already carries the column-level text-number-format style, but
# re-assert it explicitly from an existing date cell for safety.
$ws.Range("E2").Copy()
$ws.Range("E6:E12").PasteSpecial(-4122)

# Venue column ("D") alternates between the default style and the "s=2"
# style, matching rows 7, 9 and 11 in the new block (same quirk already
# present on rows 3 and 5 of the original sheet).
$ws.Range("D3").Copy()
$ws.Range("D7").PasteSpecial(-4122)
$ws.Range("D9").PasteSpecial(-4122)
$ws.Range("D11").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row heights for the new rows match the rest of the talk rows (17pt).
$ws.Range("A6:A12").RowHeight = 17

# Dimension grows automatically with the newly populated cells; move the
# active selection to the first empty row below the new data, just like
# the saved workbook.
$ws.Range("A13").Select()
